# Insert a new weekly price record as row 421, pushing the existing
# rows 421:506 down to 422:507 (dimension grows from A1:R506 to A1:R507).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(421).Insert()

$ws.Cells.Item(421, 1).Value = 7
$ws.Cells.Item(421, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(421, 3).Value = "Ñuble"
$ws.Cells.Item(421, 4).Value = 45173
$ws.Cells.Item(421, 5).Value = 16
$ws.Cells.Item(421, 6).Value = 100112006
$ws.Cells.Item(421, 7).Value = "Repollo"
$ws.Cells.Item(421, 8).Value = "Crespo record"
$ws.Cells.Item(421, 9).Value = "Primera"
$ws.Cells.Item(421, 10).Value = 250
$ws.Cells.Item(421, 11).Value = 1000
$ws.Cells.Item(421, 12).Value = 1000
$ws.Cells.Item(421, 13).Value = 1000
$ws.Cells.Item(421, 14).Value = "`$/unidad"
$ws.Cells.Item(421, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(421, 16).Value = 1000
$ws.Cells.Item(421, 17).Value = 1
$ws.Cells.Item(421, 18).Value = "Hortaliza"
